$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1943462897526502
$ws.Range("C2").Value = 0.5830388692579506
$ws.Range("J2").Value = 0.02120141342756184
$ws.Range("P2").Value = 0.1342756183745583
$ws.Range("S2").Value = 0.06713780918727916
$ws.Range("B3").Value = 0.01775147928994083
$ws.Range("C3").Value = 0.01775147928994083
$ws.Range("P3").Value = 0.7633136094674556
$ws.Range("S3").Value = 0.2011834319526627
$ws.Range("J4").Value = 0.06521739130434782
$ws.Range("P4").Value = 0.6521739130434783
$ws.Range("S4").Value = 0.2826086956521739
$ws.Range("B6").Value = 0.07317073170731707
$ws.Range("D6").Value = 0.01951219512195122
$ws.Range("F6").Value = 0.04390243902439024
$ws.Range("J6").Value = 0.2439024390243902
$ws.Range("O6").Value = 0.01463414634146342
$ws.Range("Q6").Value = 0.1512195121951219
$ws.Range("R6").Value = 0.06829268292682927
$ws.Range("S6").Value = 0.3853658536585366
$ws.Range("B7").Value = 0.08287292817679558
$ws.Range("D7").Value = 0.02209944751381215
$ws.Range("F7").Value = 0.04419889502762431
$ws.Range("J7").Value = 0.1325966850828729
$ws.Range("Q7").Value = 0.2320441988950276
$ws.Range("R7").Value = 0.04972375690607735
$ws.Range("S7").Value = 0.4364640883977901
$ws.Range("B8").Value = 0.0695970695970696
$ws.Range("D8").Value = 0.01465201465201465
$ws.Range("F8").Value = 0.05494505494505494
$ws.Range("J8").Value = 0.1135531135531136
$ws.Range("O8").Value = 0.001831501831501832
$ws.Range("Q8").Value = 0.2087912087912088
$ws.Range("R8").Value = 0.07692307692307693
$ws.Range("S8").Value = 0.4597069597069597
$ws.Range("B9").Value = 0.0718232044198895
$ws.Range("D9").Value = 0.01104972375690608
$ws.Range("E9").Value = 0.005524861878453038
$ws.Range("F9").Value = 0.02209944751381215
$ws.Range("J9").Value = 0.1270718232044199
$ws.Range("O9").Value = 0.005524861878453038
$ws.Range("Q9").Value = 0.1878453038674033
$ws.Range("R9").Value = 0.08839779005524862
$ws.Range("S9").Value = 0.4806629834254144
$ws.Range("B10").Value = 0.1032210834553441
$ws.Range("D10").Value = 0.02196193265007321
$ws.Range("E10").Value = 0.0007320644216691069
$ws.Range("F10").Value = 0.07101024890190337
$ws.Range("J10").Value = 0.1288433382137628
$ws.Range("O10").Value = 0.005124450951683748
$ws.Range("Q10").Value = 0.1961932650073206
$ws.Range("R10").Value = 0.09956076134699854
$ws.Range("S10").Value = 0.3733528550512445
$ws.Range("G11").Value = 0.1619718309859155
$ws.Range("J11").Value = 0.1056338028169014
$ws.Range("K11").Value = 0.2359154929577465
$ws.Range("L11").Value = 0.4823943661971831
$ws.Range("S11").Value = 0.01408450704225352
$ws.Range("G12").Value = 0.7857142857142857
$ws.Range("J12").Value = 0.1428571428571428
$ws.Range("K12").Value = 0.007142857142857143
$ws.Range("L12").Value = 0.01428571428571429
$ws.Range("S12").Value = 0.05
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.2708333333333333
$ws.Range("S13").Value = 0.1041666666666667
$ws.Range("F15").Value = 0.0154639175257732
$ws.Range("H15").Value = 0.134020618556701
$ws.Range("I15").Value = 0.07216494845360824
$ws.Range("J15").Value = 0.4381443298969072
$ws.Range("K15").Value = 0.08247422680412371
$ws.Range("M15").Value = 0.01030927835051546
$ws.Range("N15").Value = 0.005154639175257732
$ws.Range("O15").Value = 0.05154639175257732
$ws.Range("S15").Value = 0.1907216494845361
$ws.Range("F16").Value = 0.0106951871657754
$ws.Range("H16").Value = 0.1978609625668449
$ws.Range("I16").Value = 0.0748663101604278
$ws.Range("J16").Value = 0.4598930481283423
$ws.Range("K16").Value = 0.09625668449197861
$ws.Range("M16").Value = 0.0427807486631016
$ws.Range("O16").Value = 0.03208556149732621
$ws.Range("S16").Value = 0.0855614973262032
$ws.Range("F17").Value = 0.02061855670103093
$ws.Range("H17").Value = 0.2206185567010309
$ws.Range("I17").Value = 0.09484536082474226
$ws.Range("J17").Value = 0.4020618556701031
$ws.Range("K17").Value = 0.06597938144329897
$ws.Range("M17").Value = 0.02061855670103093
$ws.Range("O17").Value = 0.06185567010309279
$ws.Range("S17").Value = 0.1134020618556701
$ws.Range("F18").Value = 0.02314814814814815
$ws.Range("H18").Value = 0.1898148148148148
$ws.Range("I18").Value = 0.1064814814814815
$ws.Range("J18").Value = 0.3981481481481481
$ws.Range("K18").Value = 0.05555555555555555
$ws.Range("M18").Value = 0.03703703703703703
$ws.Range("O18").Value = 0.07407407407407407
$ws.Range("S18").Value = 0.1157407407407407
$ws.Range("F19").Value = 0.009305654974946313
$ws.Range("H19").Value = 0.2412312097351468
$ws.Range("I19").Value = 0.06084466714387974
$ws.Range("J19").Value = 0.3736578382247673
$ws.Range("K19").Value = 0.09806728704366499
$ws.Range("M19").Value = 0.01431639226914817
$ws.Range("O19").Value = 0.06370794559770938
$ws.Range("S19").Value = 0.1388690050107373

Write-Host "Applied 107 cell updates"
